$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Core shift"
$ws.Range("L1").Value = "Dirt"
$ws.Range("M1").Value = "Scab"
$ws.Range("O1").Value = "M/br."
$ws.Range("P1").Value = "S/POUR"
$ws.Range("A2").Value = "1- JAN-22"
$ws.Range("B2").Value = "'19"
$ws.Range("E2").Value = 0
$ws.Range("A3").Value = "2- JAN-22"
$ws.Range("B3").Value = "'19"
$ws.Range("E3").Value = 0
$ws.Range("A4").Value = "3- JAN-22"
$ws.Range("B4").Value = "'19"
$ws.Range("E4").Value = 0
$ws.Range("A5").Value = "4- JAN-22"
$ws.Range("B5").Value = "'19"
$ws.Range("E5").Value = 0
$ws.Range("A6").Value = "5- JAN-22"
$ws.Range("B6").Value = "'19"
$ws.Range("C6").Value = 64
$ws.Range("D6").Value = 45
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.03125
$ws.Range("N6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("A7").Value = "6- JAN-22"
$ws.Range("B7").Value = "'19"
$ws.Range("E7").Value = 0
$ws.Range("A8").Value = "7- JAN-22"
$ws.Range("B8").Value = "'19"
$ws.Range("E8").Value = 0
$ws.Range("A9").Value = "8- JAN-22"
$ws.Range("B9").Value = "'19"
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = 0
$ws.Range("A10").Value = "9- JAN-22"
$ws.Range("B10").Value = "'19"
$ws.Range("A11").Value = "10- JAN-22"
$ws.Range("B11").Value = "'19"
$ws.Range("C11").Value = 159
$ws.Range("D11").Value = 124
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.006289308176100629
$ws.Range("N11").Value = 1
$ws.Range("A12").Value = "11- JAN-22"
$ws.Range("B12").Value = "'19"
$ws.Range("A13").Value = "12- JAN-22"
$ws.Range("B13").Value = "'19"
$ws.Range("A14").Value = "13- JAN-22"
$ws.Range("B14").Value = "'19"
$ws.Range("A15").Value = "14- JAN-22"
$ws.Range("B15").Value = "'19"
$ws.Range("A16").Value = "15- JAN-22"
$ws.Range("B16").Value = "'19"
$ws.Range("A17").Value = "16- JAN-22"
$ws.Range("B17").Value = "'19"
$ws.Range("A18").Value = "17- JAN-22"
$ws.Range("B18").Value = "'19"
$ws.Range("A19").Value = "18- JAN-22"
$ws.Range("B19").Value = "'19"
$ws.Range("A20").Value = "19- JAN-22"
$ws.Range("B20").Value = "'19"
$ws.Range("C20").Value = 82
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("K20").ClearContents()
$ws.Range("A21").Value = "20- JAN-22"
$ws.Range("B21").Value = "'19"
$ws.Range("A22").Value = "21- JAN-22"
$ws.Range("B22").Value = "'19"
$ws.Range("A23").Value = "22- JAN-22"
$ws.Range("B23").Value = "'19"
$ws.Range("A24").Value = "23- JAN-22"
$ws.Range("B24").Value = "'19"
$ws.Range("A25").Value = "24- JAN-22"
$ws.Range("B25").Value = "'19"
$ws.Range("A26").Value = "25- JAN-22"
$ws.Range("B26").Value = "'19"
$ws.Range("A27").Value = "26- JAN-22"
$ws.Range("B27").Value = "'19"
$ws.Range("A28").Value = "27- JAN-22"
$ws.Range("B28").Value = "'19"
$ws.Range("A29").Value = "28- JAN-22"
$ws.Range("B29").Value = "'19"
$ws.Range("C29").Value = 63
$ws.Range("A30").Value = "29- JAN-22"
$ws.Range("B30").Value = "'19"
$ws.Range("A31").Value = "30- JAN-22"
$ws.Range("B31").Value = "'19"
